# Commit: "moved stuff to ssgvip repo"
#
# Rename the two worksheets to reflect the new repo/source layout, drop the
# leftover "whole sheet" selection that had been saved on the first sheet,
# and widen column G on the course table so the longer catalog text fits.

$wb = $excel.ActiveWorkbook

$wsSource  = $wb.Worksheets.Item(1)
$wsCourses = $wb.Worksheets.Item(2)

# --- Sheet names ------------------------------------------------------
$wsSource.Name  = "Original source"
$wsCourses.Name = "CMSC-Courses"

# --- Clear the stray full-sheet selection left on "Original source" ---
# (It previously held sqref="A1:XFD1048576"; normalize it back to A1 while
# restoring the originally active sheet afterwards so the active tab does
# not change.)
$wsSource.Activate()
$wsSource.Range("A1").Select()
$wsCourses.Activate()

# --- Widen column G (catalog description) on the courses sheet --------
$wsCourses.Columns.Item(7).ColumnWidth = 67.72
